# Update the worksheet date heading and the 25 division problems/answers
# in the table (5 populated rows x 5 columns) to the new day's values.
#
# Each replacement below is scoped with MatchWholeWord/MatchCase-style
# exact text via Find.Execute(..., Replace:=wdReplaceAll) against the
# full document Range, which is safe here because every "old" string is
# unique in the document.
#
# NOTE on ordering: the new text for the cell "130÷9=14, 4" -> "585÷6=97, 3"
# happens to be identical to the *old* text of another cell
# ("585÷6=97, 3" -> "151÷7=21, 4"). To avoid the first replacement's
# freshly-written text being re-matched by the second, the
# "585÷6=97, 3" -> "151÷7=21, 4" replacement is executed FIRST, while the
# original "585÷6=97, 3" text still uniquely identifies that one cell.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Find.Execute could not locate text: $oldText"
    }
}

Replace-ExactText "2025-01-03 Friday" "2025-01-04 Saturday"

# Row 1
Replace-ExactText "987÷9=109, 6" "600÷4=150, 0"
Replace-ExactText "961÷2=480, 1" "296÷5=59, 1"
Replace-ExactText "926÷9=102, 8" "843÷9=93, 6"
Replace-ExactText "526÷2=263, 0" "334÷8=41, 6"
Replace-ExactText "491÷2=245, 1" "648÷3=216, 0"

# Row 2
Replace-ExactText "245÷6=40, 5" "423÷5=84, 3"
Replace-ExactText "761÷9=84, 5" "747÷7=106, 5"
Replace-ExactText "492÷8=61, 4" "277÷6=46, 1"
Replace-ExactText "876÷9=97, 3" "331÷9=36, 7"
Replace-ExactText "996÷7=142, 2" "565÷5=113, 0"

# Row 3
Replace-ExactText "821÷2=410, 1" "223÷2=111, 1"
Replace-ExactText "664÷5=132, 4" "597÷5=119, 2"
Replace-ExactText "495÷8=61, 7" "254÷2=127, 0"
Replace-ExactText "437÷2=218, 1" "944÷6=157, 2"
# (see ordering note above) must run before "130÷9=14, 4" below
Replace-ExactText "585÷6=97, 3" "151÷7=21, 4"
Replace-ExactText "130÷9=14, 4" "585÷6=97, 3"

# Row 4
Replace-ExactText "694÷2=347, 0" "421÷3=140, 1"
Replace-ExactText "761÷4=190, 1" "492÷4=123, 0"
Replace-ExactText "888÷8=111, 0" "492÷7=70, 2"
Replace-ExactText "782÷2=391, 0" "185÷6=30, 5"

# Row 5
Replace-ExactText "624÷3=208, 0" "969÷5=193, 4"
Replace-ExactText "238÷2=119, 0" "796÷2=398, 0"
Replace-ExactText "615÷2=307, 1" "481÷8=60, 1"
Replace-ExactText "676÷2=338, 0" "332÷8=41, 4"
Replace-ExactText "603÷9=67, 0" "684÷5=136, 4"
